# B6-PowerPoint.pptx - Apr 2 2020 edit
#
# 1) Re-colour the presentation's theme (Design > Colors) from the
#    "Integral" (Red Violet) palette over to the plain default "Office"
#    palette.
# 2) Re-apply a (built-in gallery) table style to every table in the
#    deck, replacing the old custom "Table_0" style.

$p = $ppt.ActivePresentation

function Convert-HexToOleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1) Theme colours -> default "Office" palette -------------------------
# ThemeColorScheme slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$cs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Item($i).RGB = Convert-HexToOleColor $officeColors[$i - 1]
}

# --- 2) Table styles -> built-in gallery style -----------------------------
$newTableStyleId = "{4E7AD843-4BDD-4297-A7B3-D7326B537E36}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyleId)
        }
    }
}
